# Remove the trailing "Ver no Jupiter ..." and "© 2020 ..." paragraphs
# (plus the blank paragraph that separated them from the bibliography),
# as scraped at a later site build. The bibliography entry itself and the
# blank paragraph that remains before the page break are left untouched.

$d = $word.ActiveDocument

$start = $d.Paragraphs.Item(45).Range.Start
$end   = $d.Paragraphs.Item(47).Range.End

$r = $d.Range($start, $end)
$r.Delete()
